$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for the season record columns (copy formatting from an
# existing header cell so they match the rest of row 1)
$ws.Range("A1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122)

$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Season record: 90 wins, 72 losses, 0 ties for every player row
for ($row = 2; $row -le 44; $row++) {
    $ws.Cells.Item($row, 29).Value = 90
    $ws.Cells.Item($row, 30).Value = 72
    $ws.Cells.Item($row, 31).Value = 0
}
